$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.74
$ws.Range("B8").Value = 6.339
$ws.Range("B10").Value = 5.885
$ws.Range("B12").Value = 5.356999999999999
$ws.Range("C12").Value = -11.659
$ws.Range("D12").Value = -7.629
$ws.Range("D13").Value = -7.772999999999999
$ws.Range("C15").Value = -13.174
$ws.Range("C17").Value = -13.271
$ws.Range("B18").Value = 5.531000000000001
$ws.Range("D21").Value = -7.931
$ws.Range("D25").Value = -7.722
$ws.Range("C26").Value = -12.545
$ws.Range("C27").Value = -12.661
$ws.Range("C28").Value = -12.357
$ws.Range("D32").Value = -6.997999999999999
$ws.Range("D36").Value = -7.346000000000001
$ws.Range("B37").Value = 8.398
$ws.Range("C37").Value = -11.862
$ws.Range("D38").Value = -7.850999999999999
$ws.Range("D41").Value = -8.275
$ws.Range("C47").Value = -12.545
$ws.Range("D52").Value = -8.083
$ws.Range("B55").Value = 5.129
$ws.Range("D59").Value = -7.87
$ws.Range("C65").Value = -12.161
$ws.Range("D67").Value = -7.696000000000001
$ws.Range("B68").Value = 4.746
$ws.Range("C73").Value = -12.177
$ws.Range("B77").Value = 6.436
$ws.Range("B78").Value = 8.821999999999999
$ws.Range("B81").Value = 5.362
$ws.Range("B82").Value = 5.609999999999999
$ws.Range("C84").Value = -12.961
$ws.Range("D84").Value = -7.693000000000001
$ws.Range("C85").Value = -12.547
$ws.Range("D88").Value = -8.153
$ws.Range("D89").Value = -8.132999999999999
$ws.Range("C93").Value = -11.468
$ws.Range("C95").Value = -12.091
$ws.Range("D95").Value = -7.736
$ws.Range("C98").Value = -12.959
$ws.Range("C99").Value = -11.784
$ws.Range("C101").Value = -12.613
$ws.Range("D105").Value = -7.834000000000001
